# Adding 4.2 Visual Design
# (Rewriting the activity list / project plan data for the Gantt chart.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# ---- Project title (B1) ----
$ws.Range("B1").Value = "NutriMate Project"

# ---- Activity table: Name (B), Plan Start (C), Plan Duration (D) ----
$ws.Range("B5").Value = "1.1 Project Background and Overview"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1

$ws.Range("B6").Value = "1.2 Assign Project manager"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 1

$ws.Range("B7").Value = "1.3 Develop Project manager"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 2

$ws.Range("B8").Value = "1.4 Meet stakeholders"
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 3

$ws.Range("B9").Value = "2.1 Collect Requirement"
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 2

$ws.Range("B10").Value = "2.2 Create Use Case Diagram"
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = 1

$ws.Range("B11").Value = "2.3 Establish Scope Plan"
$ws.Range("C11").Value = 11
$ws.Range("D11").Value = 3

$ws.Range("B12").Value = "2.4 Create WBS"
$ws.Range("C12").Value = 14
$ws.Range("D12").Value = 1

$ws.Range("B13").Value = "2.5 Create Gantt Chart"
$ws.Range("C13").Value = 14
$ws.Range("D13").Value = 2

$ws.Range("B14").Value = "2.6 Software Design & System Components"
$ws.Range("C14").Value = 16
$ws.Range("D14").Value = 3

$ws.Range("B15").Value = "2.7 Design User Interface"
$ws.Range("C15").Value = 19
$ws.Range("D15").Value = 3

$ws.Range("B16").Value = "3.1 Develop Front-End Interface"
$ws.Range("C16").Value = 22
$ws.Range("D16").Value = 3

$ws.Range("B17").Value = "3.2 Develop Front-End Interface"
$ws.Range("C17").Value = 25
$ws.Range("D17").Value = 3

$ws.Range("B18").Value = "3.3 Integrate Features"
$ws.Range("C18").Value = 28
$ws.Range("D18").Value = 2

$ws.Range("B19").Value = "3.4 Quality Assurance & Testing"
$ws.Range("C19").Value = 30
$ws.Range("D19").Value = 3

$ws.Range("B20").Value = "3.5 Deployment"
$ws.Range("C20").Value = 33
$ws.Range("D20").Value = 3

$ws.Range("B21").Value = "4.1 Monitor Project Work"
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 2

$ws.Range("B22").Value = "4.2 Control Scope"
$ws.Range("C22").Value = 37
$ws.Range("D22").Value = 1

$ws.Range("B23").Value = "4.3 Control Schedule"
$ws.Range("C23").Value = 38
$ws.Range("D23").Value = 2

$ws.Range("B24").Value = "4.4 Perform Quality Control"
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 1

$ws.Range("B25").Value = "4.5 Report Performance"
$ws.Range("C25").Value = 40
$ws.Range("D25").Value = 2

$ws.Range("B26").Value = "5.1 Final Performance Review"
$ws.Range("C26").Value = 41
$ws.Range("D26").Value = 1

# Row 27 was previously blank; it now becomes a new activity row.
$ws.Range("B27").Value = "5.2 Prepare Final Report"
$ws.Range("C27").Value = 42
$ws.Range("D27").Value = 2

# ---- Sheet view: zoom + selection ----
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 68
$ws.Range("B2:F2").Select() | Out-Null
